$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G, pushing the old G (wait/2000) to H.
# Column F's formatting (style + width) is carried into the new column G.
$ws.Columns.Item(7).Insert()

# New header cell G1 mirrors F1 ("executeAsyncScript").
$ws.Range("G1").Value = "executeAsyncScript"

# Row 2: F2 becomes the new script body, G2 gets the new JSON target payload.
$ws.Range("G2").Value = "{""target"":""alert('Hello');""}"
$ws.Range("F2").Value = "var a=100;"

# Update the active selection to match the saved view state.
$ws.Range("F5").Select()
